$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 25 --------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A25").Value = "Demo inplannen"
$logs.Range("B25").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C25").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D25").Value = "INTERN – Intern verzoek / Actie voor medewerker"
$logs.Range("E25").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F25").Value = "2025-08-14 21:22:16"
$logs.Range("G25").Value = "Nee"
$logs.Range("H25").Value = "Ja"
$logs.Range("I25").Value = "Nee"
$logs.Range("J25").Value = "Nee"

# Extend the conditional-formatting applied ranges to cover the new row.
$dCond = $logs.Range("D2:D24").FormatConditions.Item(1)
$dCond.ModifyAppliesToRange($logs.Range("D2:D25"))

$gCond = $logs.Range("G2:G24").FormatConditions.Item(1)
$gCond.ModifyAppliesToRange($logs.Range("G2:G25"))

$hCond = $logs.Range("H2:H24").FormatConditions.Item(1)
$hCond.ModifyAppliesToRange($logs.Range("H2:H25"))

$iCond = $logs.Range("I2:I24").FormatConditions.Item(1)
$iCond.ModifyAppliesToRange($logs.Range("I2:I25"))

$jCond = $logs.Range("J2:J24").FormatConditions.Item(1)
$jCond.ModifyAppliesToRange($logs.Range("J2:J25"))

# --- Dashboard sheet: append row 7 (new category total) -------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "INTERN – Intern verzoek / Actie voor medewerker"
$dash.Range("B7").Value = 1

# --- Chart1: extend series ranges to include the new Dashboard row --------
$co = $dash.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
